$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers (e.g. "4.06") must be forced to
# Text format first, otherwise Excel auto-converts them to numeric cells - the source
# data file stores every value in these columns as text.
$textForcedCells = @("D5", "D6", "D10", "D11", "D13", "D16", "D19", "D20", "D22", "D24", "D25", "D27", "D29", "D35", "D36", "D40", "D43", "D45", "D47", "D48", "D49", "D50")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin prices / volumes / names / links.
$ws.Range("D2").Value = "26.320.85"
$ws.Range("D3").Value = "1.588.16"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D5").Value = "210.06"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.06"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.577.02"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "64.27"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.319.49"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "211.14"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "144.89"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "15.24"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "1.304.84"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.611"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  -12.23%  "
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("D43").Value = "0.767"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "62.39"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "1.723.56"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "87.74"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0982"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("E51").Value = "  -0.38%  "
